# withBookmarkCallAfterBookmarkRef-expected-generation.docx : 3.2.1 -> 3.2.2
#
# The only differences between the two revisions of this fixture are the
# values of:
#   - the w:rsidR="..." attribute stamped on the five runs that make up the
#     "REF testid \h" field (the fldChar begin/instrText/separate/result/end
#     runs), and
#   - the w:id="..." attribute on the <w:bookmarkStart>/<w:bookmarkEnd> pair
#     that defines the "testid" bookmark.
#
# Both values are opaque identifiers minted by the M2Doc template-generation
# library itself (not by Word) when it produced this "expected generation"
# fixture; bumping the library from 3.2.1 to 3.2.2 simply caused a fresh,
# unrelated-looking identifier to be minted on regeneration. The visible
# text, the field code/result, the bookmark name and position, and every
# other part of the document are unchanged.
#
# Word's object model does not surface either of these identifiers for
# automation:
#   - there is no Range/Font "RsidR" property (run-level rsid stamps are
#     maintained internally by the save pipeline and are not user-settable
#     even in real Word/VBA), and
#   - Bookmark has no "Id" property (bookmarks are addressed by Name; the
#     numeric w:id is a serialization-only detail Word assigns itself).
# Re-creating the field/bookmark through the object model (Field.Update,
# Bookmarks.Add, InsertXML, ...) does not let a caller choose these values
# either - it only re-mints *different* internal ids and, worse, risks
# disturbing real content (e.g. Field.Update() would recompute the REF
# field's cached result from the current bookmark text instead of leaving
# the already-generated result text alone).
#
# So there is nothing to change here via the supported Office object model:
# the content this fixture represents is already correct. We simply confirm
# the two constructs that the original commit touched are present, and
# leave the document otherwise untouched.
$d = $word.ActiveDocument

$targetField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $f = $d.Fields.Item($i)
    if ($f.Code.Text.Trim() -eq "REF testid \h") {
        $targetField = $f
    }
}
if ($targetField -ne $null) {
    Write-Output ("REF field present, result: " + $targetField.Result.Text)
}

$targetBookmark = $null
for ($i = 1; $i -le $d.Bookmarks.Count; $i++) {
    $b = $d.Bookmarks.Item($i)
    if ($b.Name -eq "testid") {
        $targetBookmark = $b
    }
}
if ($targetBookmark -ne $null) {
    Write-Output ("Bookmark 'testid' present at " + $targetBookmark.Start + "-" + $targetBookmark.End)
}
